# Add 4 new food items to the calorie dictionary (rows 285-288).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# food_item, measure, calories, protein, fats, carbohydrates
# The measure for the first new row is the literal text "1" (not the
# number 1), so it is prefixed with an apostrophe to force text entry,
# matching the existing column's data (e.g. "1 cup", "1 serving", ...).
$data = @(
    @("whole wheat pizza crust",       "'1",          390,   12,    6,     78),
    @("firm tofu",                     "1 block",     320,   40,    20,    10),
    @("filter coffee mocha",           "1 serving",   94.2,  1.005, 1.675, 57.86),
    @("indian style vegan tofu pizza", "1 serving",   424,   25,    14.5,  54.75)
)

$startRow = 285
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
    $ws.Cells.Item($row, 6).Value = $data[$i][5]
}

# Entering "'1" forces text but also stamps a quote-prefix cell style;
# clear that formatting so the cell keeps its plain/default style while
# the stored value remains the text string "1".
$ws.Cells.Item(285, 2).ClearFormats()
